{"js": "// Replace the arithmetic-problem text in each cell of the worksheet table.\n// The mapping below lists, in document order (row-major, top-to-bottom /\n// left-to-right), the new expression that replaces whatever text currently\n// occupies that cell. Some source expressions repeat (e.g. \"2+52=\" appears\n// twice but is replaced by two different values), so cells must be updated\n// positionally rather than via a global text search/replace.\nconst newGrid = [\n  [\"88-59=\", \"78-32=\", \"82-65=\", \"4+35=\", \"15+23=\"],\n  [\"1+62=\", \"71+10=\", \"98-7=\", \"87+5=\", \"80+18=\"],\n  [\"44+8=\", \"41-13=\", \"19-7=\", \"4+52=\", \"37+17=\"],\n  [\"75-21=\", \"23+16=\", \"26+38=\", \"15+32=\", \"24+66=\"],\n  [\"87-61=\", \"67-22=\", \"99-51=\", \"15+41=\", \"62-45=\"],\n  [\"55+40=\", \"1+80=\", \"32-25=\", \"19+25=\", \"14-4=\"],\n  [\"33+17=\", \"0+41=\", \"33+17=\", \"15+29=\", \"55-48=\"],\n  [\"76-75=\", \"67-63=\", \"14+50=\", \"86+4=\", \"80-38=\"],\n  [\"66-57=\", \"33+16=\", \"87-60=\", \"86-28=\", \"24+66=\"],\n  [\"27+4=\", \"14+51=\", \"19+11=\", \"96-11=\", \"57-23=\"],\n  [\"25+27=\", \"86-4=\", \"84-3=\", \"72+27=\", \"81-14=\"],\n  [\"7+37=\", \"50-18=\", \"44+32=\", \"45-37=\", \"74-52=\"],\n  [\"87-9=\", \"54-49=\", \"80-30=\", \"13+52=\", \"4+36=\"],\n  [\"16+69=\", \"59+14=\", \"56+21=\", \"9+9=\", \"12+67=\"],\n  [\"53-24=\", \"51-8=\", \"49-45=\", \"83-24=\", \"78-28=\"],\n  [\"9-4=\", \"99-52=\", \"77-45=\", \"96-9=\", \"41+36=\"],\n  [\"36+6=\", \"64-50=\", \"31+66=\", \"68-45=\", \"30-1=\"],\n  [\"79-20=\", \"57-20=\", \"40-17=\", \"94-40=\", \"93-44=\"],\n  [\"47+27=\", \"72-6=\", \"99-93=\", \"13+31=\", \"3+61=\"],\n  [\"67-17=\", \"94-54=\", \"34+46=\", \"81+13=\", \"40+6=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.rowCount !== newGrid.length) {\n  throw new Error(\n    `Unexpected row count: expected ${newGrid.length}, found ${table.rowCount}.`\n  );\n}\n\n// Setting Table.values rewrites each cell's text run in place while\n// preserving the existing paragraph/run formatting (font, size, alignment).\ntable.values = newGrid;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New grid of arithmetic expressions, in document (row-major) order.\n# Positional assignment is required because some source expressions repeat\n# (e.g. \"2+52=\") yet map to different replacement values depending on cell.\n$grid = @(\n  @(\"88-59=\", \"78-32=\", \"82-65=\", \"4+35=\", \"15+23=\"),\n  @(\"1+62=\", \"71+10=\", \"98-7=\", \"87+5=\", \"80+18=\"),\n  @(\"44+8=\", \"41-13=\", \"19-7=\", \"4+52=\", \"37+17=\"),\n  @(\"75-21=\", \"23+16=\", \"26+38=\", \"15+32=\", \"24+66=\"),\n  @(\"87-61=\", \"67-22=\", \"99-51=\", \"15+41=\", \"62-45=\"),\n  @(\"55+40=\", \"1+80=\", \"32-25=\", \"19+25=\", \"14-4=\"),\n  @(\"33+17=\", \"0+41=\", \"33+17=\", \"15+29=\", \"55-48=\"),\n  @(\"76-75=\", \"67-63=\", \"14+50=\", \"86+4=\", \"80-38=\"),\n  @(\"66-57=\", \"33+16=\", \"87-60=\", \"86-28=\", \"24+66=\"),\n  @(\"27+4=\", \"14+51=\", \"19+11=\", \"96-11=\", \"57-23=\"),\n  @(\"25+27=\", \"86-4=\", \"84-3=\", \"72+27=\", \"81-14=\"),\n  @(\"7+37=\", \"50-18=\", \"44+32=\", \"45-37=\", \"74-52=\"),\n  @(\"87-9=\", \"54-49=\", \"80-30=\", \"13+52=\", \"4+36=\"),\n  @(\"16+69=\", \"59+14=\", \"56+21=\", \"9+9=\", \"12+67=\"),\n  @(\"53-24=\", \"51-8=\", \"49-45=\", \"83-24=\", \"78-28=\"),\n  @(\"9-4=\", \"99-52=\", \"77-45=\", \"96-9=\", \"41+36=\"),\n  @(\"36+6=\", \"64-50=\", \"31+66=\", \"68-45=\", \"30-1=\"),\n  @(\"79-20=\", \"57-20=\", \"40-17=\", \"94-40=\", \"93-44=\"),\n  @(\"47+27=\", \"72-6=\", \"99-93=\", \"13+31=\", \"3+61=\"),\n  @(\"67-17=\", \"94-54=\", \"34+46=\", \"81+13=\", \"40+6=\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif ($rowCount -ne $grid.Count) {\n  throw \"Unexpected row count: expected $($grid.Count), found $rowCount.\"\n}\nif ($colCount -ne $grid[0].Count) {\n  throw \"Unexpected column count: expected $($grid[0].Count), found $colCount.\"\n}\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $t.Cell($r, $c).Range.Text = $grid[$r - 1][$c - 1]\n  }\n}\n\nWrite-Output \"updated $rowCount x $colCount cells\"\n"}
